$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'68.207.02"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "'  -3.61%  "
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(3, 4).Value = "'3.690.29"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "'  -4.38%  "
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "'  -0.18%  "
$ws.Cells.Item(4, 5).Style = "Normal"
$ws.Cells.Item(5, 4).Value = "'593.84"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "'  +0.35%  "
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Cells.Item(6, 4).Value = "'181.70"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "'  +8.97%  "
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Cells.Item(7, 4).Value = "'3.677.86"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "'  -4.60%  "
$ws.Cells.Item(7, 5).Style = "Normal"
$ws.Cells.Item(8, 4).Value = "'0.626"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "'  -6.58%  "
$ws.Cells.Item(8, 5).Style = "Normal"
$ws.Cells.Item(9, 4).Value = "'1.00"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "'  -0.01%  "
$ws.Cells.Item(9, 5).Style = "Normal"
$ws.Cells.Item(10, 4).Value = "'0.714"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "'  -4.63%  "
$ws.Cells.Item(10, 5).Style = "Normal"
$ws.Cells.Item(11, 4).Value = "'0.162"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "'  -7.11%  "
$ws.Cells.Item(11, 5).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "'  +4.69%  "
$ws.Cells.Item(12, 5).Style = "Normal"
$ws.Cells.Item(13, 4).Value = "'0.0000290"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "'  -9.35%  "
$ws.Cells.Item(13, 5).Style = "Normal"
$ws.Cells.Item(14, 4).Value = "'10.33"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "'  -9.62%  "
$ws.Cells.Item(14, 5).Style = "Normal"
$ws.Cells.Item(15, 4).Value = "'4.257.17"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "'  -5.12%  "
$ws.Cells.Item(15, 5).Style = "Normal"
$ws.Cells.Item(16, 4).Value = "'3.684.37"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "'  -5.12%  "
$ws.Cells.Item(16, 5).Style = "Normal"
$ws.Cells.Item(17, 4).Value = "'19.28"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "'  -8.58%  "
$ws.Cells.Item(17, 5).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "'  -2.34%  "
$ws.Cells.Item(18, 5).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "'  -7.02%  "
$ws.Cells.Item(19, 5).Style = "Normal"
$ws.Cells.Item(20, 4).Value = "'12.76"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "'  -7.14%  "
$ws.Cells.Item(20, 5).Style = "Normal"
$ws.Cells.Item(21, 4).Value = "'67.856.35"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "'  -4.19%  "
$ws.Cells.Item(21, 5).Style = "Normal"
$ws.Cells.Item(22, 4).Value = "'408.72"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "'  -6.21%  "
$ws.Cells.Item(22, 5).Style = "Normal"
$ws.Cells.Item(23, 4).Value = "'4.53"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "'  -3.44%  "
$ws.Cells.Item(23, 5).Style = "Normal"
$ws.Cells.Item(24, 4).Value = "'88.54"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "'  -5.98%  "
$ws.Cells.Item(24, 5).Style = "Normal"
$ws.Cells.Item(25, 4).Value = "'3.01"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "'  -7.81%  "
$ws.Cells.Item(25, 5).Style = "Normal"
$ws.Cells.Item(26, 4).Value = "'12.75"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "'  -7.88%  "
$ws.Cells.Item(26, 5).Style = "Normal"
$ws.Cells.Item(27, 4).Value = "'10.88"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "'  -3.61%  "
$ws.Cells.Item(27, 5).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "'  -6.00%  "
$ws.Cells.Item(28, 5).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "'  +2.06%  "
$ws.Cells.Item(29, 5).Style = "Normal"
$ws.Cells.Item(30, 4).Value = "'9.40"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "'  -9.09%  "
$ws.Cells.Item(30, 5).Style = "Normal"
$ws.Cells.Item(31, 4).Value = "'32.68"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "'  -6.77%  "
$ws.Cells.Item(31, 5).Style = "Normal"
$ws.Cells.Item(32, 4).Value = "'7.27"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "'  -11.40%  "
$ws.Cells.Item(32, 5).Style = "Normal"
$ws.Cells.Item(33, 4).Value = "'12.43"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "'  -8.26%  "
$ws.Cells.Item(33, 5).Style = "Normal"
$ws.Cells.Item(34, 4).Value = "'0.117"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "'  -7.00%  "
$ws.Cells.Item(34, 5).Style = "Normal"
$ws.Cells.Item(35, 4).Value = "'43.49"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "'  -10.58%  "
$ws.Cells.Item(35, 5).Style = "Normal"
$ws.Cells.Item(36, 4).Value = "'64.04"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "'  -7.76%  "
$ws.Cells.Item(36, 5).Style = "Normal"
$ws.Cells.Item(37, 4).Value = "'598.31"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "'  -4.97%  "
$ws.Cells.Item(37, 5).Style = "Normal"
$ws.Cells.Item(38, 4).Value = "'0.0₃0884"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "'  -9.99%  "
$ws.Cells.Item(38, 5).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "'  +0.10%  "
$ws.Cells.Item(39, 5).Style = "Normal"
$ws.Cells.Item(40, 4).Value = "'0.398"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "'  -6.86%  "
$ws.Cells.Item(40, 5).Style = "Normal"
$ws.Cells.Item(41, 4).Value = "'1.00"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "'  -0.21%  "
$ws.Cells.Item(41, 5).Style = "Normal"
$ws.Cells.Item(42, 4).Value = "'0.136"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "'  -7.87%  "
$ws.Cells.Item(42, 5).Style = "Normal"
$ws.Cells.Item(43, 4).Value = "'2.80"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "'  +3.13%  "
$ws.Cells.Item(43, 5).Style = "Normal"
$ws.Cells.Item(44, 4).Value = "'2.99"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "'  -8.73%  "
$ws.Cells.Item(44, 5).Style = "Normal"
$ws.Cells.Item(45, 4).Value = "'0.0435"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "'  -7.13%  "
$ws.Cells.Item(45, 5).Style = "Normal"
$ws.Cells.Item(46, 4).Value = "'2.87"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "'  -10.95%  "
$ws.Cells.Item(46, 5).Style = "Normal"
$ws.Cells.Item(47, 4).Value = "'9.16"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "'  -8.99%  "
$ws.Cells.Item(47, 5).Style = "Normal"
$ws.Cells.Item(48, 4).Value = "'2.72"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "'  -3.92%  "
$ws.Cells.Item(48, 5).Style = "Normal"
$ws.Cells.Item(49, 2).Value = "ApeXProtocol"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(49, 4).Value = "'3.18"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "'  -5.89%  "
$ws.Cells.Item(49, 5).Style = "Normal"
$ws.Cells.Item(50, 2).Value = "Stellar"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(50, 4).Value = "'0.134"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "'  -6.81%  "
$ws.Cells.Item(50, 5).Style = "Normal"
$ws.Cells.Item(51, 2).Value = "Maker"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(51, 4).Value = "'2.733.24"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "'  -3.65%  "
$ws.Cells.Item(51, 5).Style = "Normal"

Write-Host "Applied all crypto list updates"